$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last name in row 2 from "Jones" to "Jane"
$ws.Range("B2").Value = "Jane"

# Clear the Job_ID value in row 5 and replace with a single space
$ws.Range("C5").Value = " "

# Remove the entire "Salary" column (column D)
$ws.Range("D1:D5").Delete()

# Update the active selection to F5
$ws.Range("F5").Select()
